$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.1550416982562547,
    0.5504169825625473,
    0.8699772554965883,
    0.9571645185746778,
    0.9859742228961335,
    0.9950720242608037,
    0.9973464746019712,
    0.9977255496588325,
    0.9977255496588325,
    0.9981046247156938,
    0.9981046247156938,
    0.9984836997725551,
    0.9988627748294164,
    0.9988627748294164,
    0.9988627748294164,
    0.9988627748294164,
    0.9988627748294164,
    0.9988627748294164,
    0.9992418498862777,
    0.9992418498862777,
    0.9992418498862777,
    0.9992418498862777,
    0.9992418498862777,
    0.9992418498862777,
    0.999620924943139,
    0.999620924943139,
    0.999620924943139,
    0.999620924943139,
    0.999620924943139
)

$row = 2
foreach ($v in $values) {
    $ws.Cells.Item($row, 1).Value = $v
    $row = $row + 1
}
